# Stundenliste-SebastianEderer.xlsx update
# - added custom zooming to apexchart
#
# Fills in the hours-worked table for the week of 24.02.2021-03.03.2021
# (rows 24-30), records the two trailing dates (rows 31-32) and drops the
# now-unused tail rows (33-35) that only held the carried-forward total.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: stamp a date serial into $addr while re-using the same "short
# date" cell style already applied to the existing Datum column (A23),
# so we don't spawn a brand-new custom number format in styles.xml.
function Set-DateCell {
    param($addr, $serial)
    $ws.Range("A23").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range($addr).Value = $serial
}

# --- Datum column first (matches the shared-string allocation order of
#     the authored workbook: the literal text "29.02.2021" in A29 is
#     registered before the Tätigkeit strings that follow it) ---
Set-DateCell "A24" 44251
$ws.Range("B24").Value = 5

Set-DateCell "A25" 44252
$ws.Range("B25").Value = 7

Set-DateCell "A26" 44253
$ws.Range("B26").Value = 8

Set-DateCell "A27" 44254
$ws.Range("B27").Value = 8

Set-DateCell "A28" 44255
$ws.Range("B28").Value = 7

# Row 29's Datum cell was typed as literal text, not a real date.
$ws.Range("A29").Value = "29.02.2021"
$ws.Range("B29").Value = 9

Set-DateCell "A30" 44256
$ws.Range("B30").Value = 8

# --- Tätigkeit column, row 24 already set inline above for ordering;
#     remaining rows filled now ---
$ws.Range("D24").Value = "Poster + Development WebUI"
$ws.Range("D25").Value = "Research on alternatives for apexcharts"
$ws.Range("D26").Value = "Development WebUI"
$ws.Range("D27").Value = "Meeting + Development WebUI"
$ws.Range("D28").Value = "Development WebUI"
$ws.Range("D29").Value = "Research on custom range sliders + Development WebUI"
$ws.Range("D30").Value = "Development WebUI"

# --- Two more dates with no hours/activity recorded yet ---
Set-DateCell "A31" 44257
Set-DateCell "A32" 44258

# These stale rows only ever held the running-total formula; clear it
# out of 31/32 and drop 33-35 completely now that the data stops at 32.
$ws.Range("C31:C32").ClearContents()
$ws.Range("33:35").Delete()

# --- Viewport: scrolled down with I28 as the active (if empty) cell ---
$ws.Range("I28").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
